$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.822.89"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.289.77"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "2.646.58"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "2.293.11"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.772"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "42.751.47"
$ws.Range("E19").Value = "  -5.47%  "
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.21%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0683"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("E40").Value = "  -4.54%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").Value = "2.013.07"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").Value = "2.514.28"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.83%  "
